# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly generated output numbers.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 42
$ws1.Range("F5").Value = 2777
$ws1.Range("F6").Value = 1942
$ws1.Range("F9").Value = 991
$ws1.Range("F10").Value = 187
$ws1.Range("F11").Value = 28

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 42
$ws4.Range("F5").Value = 2777
$ws4.Range("F6").Value = 1942
$ws4.Range("F10").Value = 991
$ws4.Range("F11").Value = 187
$ws4.Range("F12").Value = 28

$wb.Save()
